$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cell updates (Coin name / Link / Volume columns) ---
$textUpdates = @{
    'E2' = '  -1.50%  '
    'E3' = '  -4.83%  '
    'E4' = '  -0.07%  '
    'E5' = '  -2.47%  '
    'E6' = '  -0.33%  '
    'E7' = '  +0.01%  '
    'E8' = '  +21.07%  '
    'E9' = '  -5.16%  '
    'E10' = '  -0.93%  '
    'E11' = '  -6.30%  '
    'E12' = '  -1.49%  '
    'E13' = '  +1.24%  '
    'E15' = '  -1.60%  '
    'E16' = '  -5.10%  '
    'E17' = '  -4.70%  '
    'E18' = '  -5.23%  '
    'E19' = '  +1.09%  '
    'E20' = '  -0.92%  '
    'E21' = '  -4.67%  '
    'E22' = '  -0.04%  '
    'E23' = '  -2.48%  '
    'E24' = '  -3.95%  '
    'E25' = '  +0.12%  '
    'E26' = '  -4.25%  '
    'E27' = '  -5.74%  '
    'E28' = '  -5.66%  '
    'E29' = '  -8.13%  '
    'E30' = '  +0.17%  '
    'E31' = '  -5.48%  '
    'E32' = '  -1.51%  '
    'E33' = '  -0.53%  '
    'E34' = '  -2.18%  '
    'E35' = '  -1.86%  '
    'B36' = 'NEARProtocol'
    'C36' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E36' = '  -3.13%  '
    'B37' = 'ImmutableX'
    'C37' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'E37' = '  -3.95%  '
    'E38' = '  -5.93%  '
    'E39' = '  +10.21%  '
    'E40' = '  -0.55%  '
    'E41' = '  +0.16%  '
    'E42' = '  -0.11%  '
    'E43' = '  -3.92%  '
    'E44' = '  -4.54%  '
    'E45' = '  -6.80%  '
    'E46' = '  -0.56%  '
    'E47' = '  -4.23%  '
    'E48' = '  -3.01%  '
    'E49' = '  -7.16%  '
    'B50' = 'Maker'
    'C50' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'E50' = '  -4.43%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E51' = '  -5.38%  '
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Price column (D) updates ---
# These look numeric (e.g. "1.00", "474.97"), so a leading apostrophe
# is used to force Excel to store them as text, same as the original
# inline-string cells. The apostrophe itself is not stored as part of
# the value, it only sets the "quote prefix" flag on the cell, which
# we immediately clear below by reapplying the default "Normal" style
# so the resulting cell matches a plain, unstyled text cell.
$priceUpdates = @{
    'D2' = '''55.150.61'
    'D3' = '''2.342.02'
    'D4' = '''1.00'
    'D5' = '''474.97'
    'D6' = '''144.85'
    'D8' = '''0.614'
    'D9' = '''2.338.82'
    'D10' = '''0.0958'
    'D11' = '''5.43'
    'D12' = '''0.326'
    'D14' = '''2.747.70'
    'D15' = '''55.118.01'
    'D16' = '''19.94'
    'D18' = '''2.346.77'
    'D19' = '''4.55'
    'D20' = '''313.57'
    'D21' = '''9.55'
    'D23' = '''5.63'
    'D24' = '''56.01'
    'D26' = '''0.393'
    'D27' = '''0.151'
    'D28' = '''2.434.48'
    'D29' = '''7.04'
    'D31' = '''0.0₃0738'
    'D33' = '''18.10'
    'D34' = '''1.47'
    'D35' = '''5.07'
    'D36' = '''3.60'
    'D37' = '''1.09'
    'D38' = '''0.808'
    'D40' = '''33.60'
    'D41' = '''0.998'
    'D43' = '''3.37'
    'D44' = '''0.575'
    'D45' = '''0.0515'
    'D47' = '''248.70'
    'D48' = '''0.0220'
    'D49' = '''4.35'
    'D50' = '''1.784.88'
    'D51' = '''16.58'
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = $priceUpdates[$addr]
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Style = "Normal"
}
